$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every existing row down by one to make room for the new "Date and Time" row
$ws.Rows.Item(1).Insert()

$ws.Cells.Item(1,1).Value = "Date and Time"
$ws.Cells.Item(1,2).Value = "2024-03-12 19:35:27.687000 to 2024-03-12 20:40:53.242000"
$ws.Cells.Item(2,1).Value = "Total time taken for the ride"
$ws.Cells.Item(2,2).Value = 0.04510481481481481
$ws.Cells.Item(2,2).NumberFormat = "[hh]:mm:ss"
$ws.Cells.Item(3,1).Value = "Actual Ampere-hours (Ah)"
$ws.Cells.Item(3,2).Value = 27.0011225
$ws.Cells.Item(4,1).Value = "Actual Watt-hours (Wh)"
$ws.Cells.Item(4,2).Value = 1384.736395011389
$ws.Cells.Item(5,1).Value = "Starting SoC (Ah)"
$ws.Cells.Item(5,2).Value = 39.556
$ws.Cells.Item(6,1).Value = "Ending SoC (Ah)"
$ws.Cells.Item(6,2).Value = 10.745
$ws.Cells.Item(7,1).Value = "Starting SoC (%)"
$ws.Cells.Item(7,2).Value = 27
$ws.Cells.Item(8,1).Value = "Ending SoC (%)"
$ws.Cells.Item(8,2).Value = 99
$ws.Cells.Item(9,1).Value = "Total distance covered (km)"
$ws.Cells.Item(9,2).Value = 40.59608462201849
$ws.Cells.Item(10,1).Value = "Total energy consumption(WH/KM)"
$ws.Cells.Item(10,2).Value = 34.11009726441294
$ws.Cells.Item(11,1).Value = "Total SOC consumed(%)"
$ws.Cells.Item(11,2).Value = 72
$ws.Cells.Item(12,1).Value = "Mode"
$ws.Cells.Item(12,2).Value = "Custom mode`n85.39%`nSports mode`n9.95%`nEco mode`n2.11%"
$ws.Cells.Item(13,1).Value = "Peak Power(kW)"
$ws.Cells.Item(13,2).Value = 6239.258608
$ws.Cells.Item(14,1).Value = "Average Power(kW)"
$ws.Cells.Item(14,2).Value = -1290.795189549715
$ws.Cells.Item(15,1).Value = "Total Energy Regenerated(kWh)"
$ws.Cells.Item(15,2).Value = 104.6671430580555
$ws.Cells.Item(16,1).Value = "Regenerative Effectiveness(%)"
$ws.Cells.Item(16,2).Value = 7.027453633803263
$ws.Cells.Item(17,1).Value = "Highest Cell Voltage(V)"
$ws.Cells.Item(17,2).Value = 3.466
$ws.Cells.Item(18,1).Value = "Lowest Cell Voltage(V)"
$ws.Cells.Item(18,2).Value = 3.077
$ws.Cells.Item(19,1).Value = "Difference in Cell Voltage(V)"
$ws.Cells.Item(19,2).Value = 0.3890000000000002
$ws.Cells.Item(20,1).Value = "Minimum Temperature(C)"
$ws.Cells.Item(20,2).Value = 35
$ws.Cells.Item(21,1).Value = "Maximum Temperature(C)"
$ws.Cells.Item(21,2).Value = 47
$ws.Cells.Item(22,1).Value = "Difference in Temperature(C)"
$ws.Cells.Item(22,2).Value = 12
$ws.Cells.Item(23,1).Value = "Maximum Fet Temperature-BMS(C)"
$ws.Cells.Item(23,2).Value = 71
$ws.Cells.Item(24,1).Value = "Maximum Afe Temperature-BMS(C)"
$ws.Cells.Item(24,2).Value = 67
$ws.Cells.Item(25,1).Value = "Maximum PCB Temperature-BMS(C)"
$ws.Cells.Item(25,2).Value = 67
$ws.Cells.Item(26,1).Value = "Maximum MCU Temperature(C)"
$ws.Cells.Item(26,2).Value = 47
$ws.Cells.Item(27,1).Value = "Maximum Motor Temperature(C)"
$ws.Cells.Item(27,2).Value = 0
$ws.Cells.Item(28,1).Value = "Abnormal Motor Temperature Detected(C)"
$ws.Cells.Item(28,2).Value = 0
$ws.Cells.Item(29,1).Value = "highest cell temp(C)"
$ws.Cells.Item(29,2).Value = 47
$ws.Cells.Item(30,1).Value = "lowest cell temp(C)"
$ws.Cells.Item(30,2).Value = 35
$ws.Cells.Item(31,1).Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Cells.Item(31,2).Value = 12
$ws.Cells.Item(32,1).Value = "Battery Voltage(V)"
$ws.Cells.Item(32,2).Value = 56
$ws.Cells.Item(33,1).Value = "Total energy charged(kWh)"
$ws.Cells.Item(33,2).Value = 1.51206286
$ws.Cells.Item(34,1).Value = "Electricity consumption units(kW)"
$ws.Cells.Item(34,2).Value = 0.0000001077796923559433
$ws.Cells.Item(35,1).Value = "Cycle Count of battery"
$ws.Cells.Item(35,2).Value = 116
$ws.Cells.Item(36,1).Value = "Idling time percentage"
$ws.Cells.Item(36,2).Value = 4.500972235504128
$ws.Cells.Item(37,1).Value = "Time spent in 0-10 km/h"
$ws.Cells.Item(37,2).Value = 8.112588059035415
$ws.Cells.Item(38,1).Value = "Time spent in 10-20 km/h"
$ws.Cells.Item(38,2).Value = 8.380351279844442
$ws.Cells.Item(39,1).Value = "Time spent in 20-30 km/h"
$ws.Cells.Item(39,2).Value = 12.06528322335915
$ws.Cells.Item(40,1).Value = "Time spent in 30-40 km/h"
$ws.Cells.Item(40,2).Value = 17.33766854738453
$ws.Cells.Item(41,1).Value = "Time spent in 40-50 km/h"
$ws.Cells.Item(41,2).Value = 17.84769372987791
$ws.Cells.Item(42,1).Value = "Time spent in 50-60 km/h"
$ws.Cells.Item(42,2).Value = 23.9361193458927
$ws.Cells.Item(43,1).Value = "Time spent in 60-70 km/h"
$ws.Cells.Item(43,2).Value = 7.586624589589111
$ws.Cells.Item(44,1).Value = "Time spent in 70-80 km/h"
$ws.Cells.Item(44,2).Value = 0.0637531478116732
$ws.Cells.Item(45,1).Value = "Time spent in 80-90 km/h"
$ws.Cells.Item(45,2).Value = 0

$ws.Cells.Item(1,2).Style = "Normal"
